# "Generate Report for Handoff"
# Updates the localization-status report: status moves from "In Translation"
# to "Ready for handoff" and the handoff timestamps are refreshed. Columns
# that display the new (longer) status text are then resized to fit it.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-30 04:58:31"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("H2").Value = "2016-08-30 04:58:27"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("H2").Value = "2016-08-30 04:58:31"

# --- Widen the "Status" columns so the longer "Ready for handoff" text fits ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333
$wsZh.Columns.Item(3).ColumnWidth = 16.3333333333333
$wsDe.Columns.Item(3).ColumnWidth = 16.3333333333333
